# Add two new columns, I ("I0") and J ("IF"), to the single worksheet.
# Header row (row 1) gets the same style as the other header cells (s="1"),
# and every data row (2-61) gets a numeric value in I and J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Mirror the existing header formatting (bold, centered/top-aligned,
# thin box border) from column H onto the two new header cells.
$srcHeader = $ws.Range("H1")
$newHeaders = $ws.Range("I1:J1")
$newHeaders.Font.Bold = $srcHeader.Font.Bold
$newHeaders.HorizontalAlignment = $srcHeader.HorizontalAlignment
$newHeaders.VerticalAlignment = $srcHeader.VerticalAlignment
$newHeaders.Borders.LineStyle = $srcHeader.Borders.Item(7).LineStyle

# --- Data values --------------------------------------------------------
$i0 = 8,9,8,8,9,9,7,7,8,8,6,7,7,6,5,6,9,7,8,9,3,9,4,7,8,10,9,6,8,10,10,7,7,9,7,6,6,7,9,7,6,7,7,8,6,9,9,6,10,7,7,7,8,6,7,6,7,6,5,3
$iF = 8,9,8,8,9,9,7,7,8,8,7,7,7,6,6,6,9,7,9,9,4,9,5,7,8,10,9,7,8,10,10,8,7,9,7,6,6,7,9,8,6,8,7,8,6,9,9,6,10,7,7,7,8,6,7,6,7,6,5,3

for ($r = 2; $r -le 61; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0[$idx]
    $ws.Cells.Item($r, 10).Value = $iF[$idx]
}
